# aggiornamento fino a 20/09/2021
# Extend the daily series with 11 new rows (375-385), one per day from
# 2021-09-10 (serial 44449) through 2021-09-20 (serial 44459), mirroring
# the existing rows: date in column A (formatted/bordered like the prior
# row), and 0 in columns B, C, D.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 374
$firstNewRow = $lastRow + 1
$lastNewRow = 385
$firstSerial = 44449

for ($row = $firstNewRow; $row -le $lastNewRow; $row++) {
    $serial = $firstSerial + ($row - $firstNewRow)
    $ws.Cells.Item($row, 1).Value = $serial
    $ws.Cells.Item($row, 2).Value = 0
    $ws.Cells.Item($row, 3).Value = 0
    $ws.Cells.Item($row, 4).Value = 0
}

# Carry the date column's formatting (border, bold, centered, custom
# yyyy-mm-dd date format) down onto the new rows, same as the rest of
# column A.
$ws.Range("A" + $lastRow).Copy()
$ws.Range("A" + $firstNewRow + ":A" + $lastNewRow).PasteSpecial(-4122)
$excel.CutCopyMode = 0
